$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '31.018.43'
$ws.Range('E2').Value = '  +1.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.955.41'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.17'
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4867'
$ws.Range('E7').Value = '  +1.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2947'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06822'
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.16'
$ws.Range('E10').Value = '  -1.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '107.14'
$ws.Range('E11').Value = '  -3.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.954.97'
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07819'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.456'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7019'
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '282.81'
$ws.Range('E16').Value = '  -3.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '31.037.59'
$ws.Range('E17').Value = '  +1.24%  '
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007685'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.210.28'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.492'
$ws.Range('E22').Value = '  -2.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('E24').Value = '  -1.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.809'
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.11'
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.95'
$ws.Range('E27').Value = '  -1.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.198'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1054'
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.418'
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.583'
$ws.Range('E31').Value = '  -1.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.605'
$ws.Range('E32').Value = '  -2.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.439'
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04930'
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7632'
$ws.Range('E35').Value = '  -2.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.170'
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.730'
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02005'
$ws.Range('E38').Value = '  -2.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.703'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.522'
$ws.Range('E40').Value = '  +5.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.101'
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.86'
$ws.Range('E42').Value = '  +6.96%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8877'
$ws.Range('E43').Value = '  +1.56%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4460'
$ws.Range('E44').Value = '  -0.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '109.24'
$ws.Range('E45').Value = '  -1.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.207'
$ws.Range('E46').Value = '  +10.88%  '
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '997.93'
$ws.Range('E48').Value = '  +9.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1256'
$ws.Range('E49').Value = '  -1.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.299'
$ws.Range('E50').Value = '  -0.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2592'
$ws.Range('E51').Value = '  +3.36%  '
